# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets, reflecting newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 714
    7  = 21
    10 = 2
    11 = 4584
    12 = 4411
    15 = 150
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
